$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1475
$ws.Range("I18").Value = 1475
$ws.Range("K18").Value = 1475
$ws.Range("M18").Value = -1191

$ws.Range("H39").Value = 6978.9443
$ws.Range("I39").Value = 129.27272
$ws.Range("J39").Value = 17742.715
$ws.Range("K39").Value = 387.81816
$ws.Range("L39").Value = 53228.145
$ws.Range("M39").Value = -91.81815999999998
$ws.Range("N39").Value = -53820.145

$ws.Range("H40").Value = 2274.6667
$ws.Range("I40").Value = 2274.6667
$ws.Range("K40").Value = 2274.6667
$ws.Range("M40").Value = -2099.6667

$ws.Range("H43").Value = 6223.154
$ws.Range("J43").Value = 6112.5
$ws.Range("L43").Value = 6112.5
$ws.Range("N43").Value = -6250.5

$ws.Range("H61").Value = 1253.9
$ws.Range("I61").Value = 948.7778
$ws.Range("K61").Value = 2846.3334
$ws.Range("M61").Value = -2674.3334

$ws.Range("H64").Value = 5974.5
$ws.Range("I64").Value = 5800
$ws.Range("J64").Value = 6149
$ws.Range("K64").Value = 5800
$ws.Range("L64").Value = 6149
$ws.Range("M64").Value = -5552
$ws.Range("N64").Value = -6645

$ws.Range("H67").Value = 5974.5
$ws.Range("I67").Value = 5800
$ws.Range("J67").Value = 6149
$ws.Range("K67").Value = 5800
$ws.Range("L67").Value = 6149
$ws.Range("M67").Value = -4942
$ws.Range("N67").Value = -7865

$ws.Range("H74").Value = 4547.6665
$ws.Range("I74").Value = 3711.2856
$ws.Range("K74").Value = 3711.2856
$ws.Range("M74").Value = -2775.2856

$ws.Range("H77").Value = 4547.6665
$ws.Range("I77").Value = 3711.2856
$ws.Range("K77").Value = 18556.428
$ws.Range("M77").Value = -13876.428

$ws.Range("H80").Value = 545.1818
$ws.Range("J80").Value = 299.5
$ws.Range("L80").Value = 898.5
$ws.Range("N80").Value = -2894.5

$ws.Range("H83").Value = 545.1818
$ws.Range("J83").Value = 299.5
$ws.Range("L83").Value = 2695.5
$ws.Range("N83").Value = -12679.5

$ws.Range("H87").Value = 188111.73
$ws.Range("J87").Value = 196895.72
$ws.Range("L87").Value = 196895.72
$ws.Range("N87").Value = -199391.72

$ws.Range("H90").Value = 188111.73
$ws.Range("J90").Value = 196895.72
$ws.Range("L90").Value = 590687.16
$ws.Range("N90").Value = -603167.16

$ws.Range("H98").Value = 3517.8096
$ws.Range("I98").Value = 3517.8096
$ws.Range("K98").Value = 3517.8096
$ws.Range("M98").Value = -2019.8096

$ws.Range("H106").Value = 3077.2856
$ws.Range("I106").Value = 2965.1667
$ws.Range("J106").Value = 3750
$ws.Range("K106").Value = 2965.1667
$ws.Range("L106").Value = 3750
$ws.Range("M106").Value = -2334.1667
$ws.Range("N106").Value = -5012

$ws.Range("H113").Value = 12741.363
$ws.Range("I113").Value = 41249.5
$ws.Range("J113").Value = 6406.222
$ws.Range("K113").Value = 41249.5
$ws.Range("L113").Value = 6406.222
$ws.Range("M113").Value = -37995.5
$ws.Range("N113").Value = -12914.222

$ws.Range("H116").Value = 3793.2
$ws.Range("I116").Value = 3306.1428
$ws.Range("J116").Value = 4929.6665
$ws.Range("K116").Value = 3306.1428
$ws.Range("L116").Value = 4929.6665
$ws.Range("M116").Value = 135.8571999999999
$ws.Range("N116").Value = -11813.6665

$ws.Range("H122").Value = 3517.8096
$ws.Range("I122").Value = 3517.8096
$ws.Range("K122").Value = 10553.4288
$ws.Range("M122").Value = -8103.4288

$ws.Range("H129").Value = 1327.1364
$ws.Range("I129").Value = 626.5
$ws.Range("K129").Value = 1879.5
$ws.Range("M129").Value = 3120.5

$ws.Range("H132").Value = 8688.333000000001
$ws.Range("I132").Value = 2448.7896
$ws.Range("J132").Value = 32398.6
$ws.Range("K132").Value = 7346.3688
$ws.Range("L132").Value = 97195.79999999999
$ws.Range("M132").Value = -4816.3688
$ws.Range("N132").Value = -102255.8

$ws.Range("H133").Value = 50000
$ws.Range("J133").Value = 50000
$ws.Range("L133").Value = 50000
$ws.Range("N133").Value = -60120

$ws.Range("H135").Value = 7503.067
$ws.Range("I135").Value = 8388.154
$ws.Range("K135").Value = 75493.386
$ws.Range("M135").Value = -72958.386

$ws.Range("H137").Value = 19651.387
$ws.Range("I137").Value = 11866.477
$ws.Range("J137").Value = 35999.7
$ws.Range("K137").Value = 35599.431
$ws.Range("L137").Value = 107999.1
$ws.Range("M137").Value = -33049.431
$ws.Range("N137").Value = -113099.1

$ws.Range("H138").Value = 4317.421
$ws.Range("I138").Value = 6773.909
$ws.Range("J138").Value = 3316.6296
$ws.Range("K138").Value = 20321.727
$ws.Range("L138").Value = 9949.888800000001
$ws.Range("M138").Value = -15181.727
$ws.Range("N138").Value = -20229.8888

$ws.Range("H141").Value = 3497.7778
$ws.Range("I141").Value = 3497.7778
$ws.Range("K141").Value = 10493.3334
$ws.Range("M141").Value = -5313.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 9267.727999999999
$ws.Range("I2").Value = 9368.5
$ws.Range("K2").Value = 9368.5
$ws.Range("M2").Value = -9255.5

$ws.Range("H32").Value = 4646.018
$ws.Range("I32").Value = 2180.9443
$ws.Range("J32").Value = 9316.684999999999
$ws.Range("K32").Value = 2180.9443
$ws.Range("L32").Value = 9316.684999999999
$ws.Range("M32").Value = -1893.9443
$ws.Range("N32").Value = -9890.684999999999

$ws.Range("H61").Value = 690346.1
$ws.Range("I61").Value = 4095.9412
$ws.Range("K61").Value = 4095.9412
$ws.Range("M61").Value = -3883.9412

$ws.Range("H63").Value = 2398.25
$ws.Range("I63").Value = 1531.1666
$ws.Range("J63").Value = 4999.5
$ws.Range("K63").Value = 1531.1666
$ws.Range("L63").Value = 4999.5
$ws.Range("M63").Value = -845.1666
$ws.Range("N63").Value = -6371.5

$ws.Range("H66").Value = 2398.25
$ws.Range("I66").Value = 1531.1666
$ws.Range("J66").Value = 4999.5
$ws.Range("K66").Value = 7655.833000000001
$ws.Range("L66").Value = 24997.5
$ws.Range("M66").Value = -4223.833000000001
$ws.Range("N66").Value = -31861.5

$ws.Range("H74").Value = 9557.097
$ws.Range("I74").Value = 3163.1538
$ws.Range("K74").Value = 3163.1538
$ws.Range("M74").Value = -2289.1538

$ws.Range("H77").Value = 9557.097
$ws.Range("I77").Value = 3163.1538
$ws.Range("K77").Value = 15815.769
$ws.Range("M77").Value = -11447.769

$ws.Range("H88").Value = 1934.9048
$ws.Range("I88").Value = 1900.5
$ws.Range("K88").Value = 1900.5
$ws.Range("M88").Value = -1494.5

$ws.Range("H91").Value = 1934.9048
$ws.Range("I91").Value = 1900.5
$ws.Range("K91").Value = 1900.5
$ws.Range("M91").Value = -496.5

$ws.Range("H102").Value = 4930.1333
$ws.Range("I102").Value = 5069.5
$ws.Range("K102").Value = 5069.5
$ws.Range("M102").Value = -3447.5

$ws.Range("H110").Value = 7066.375
$ws.Range("I110").Value = 9089
$ws.Range("K110").Value = 9089
$ws.Range("M110").Value = -7044

$ws.Range("H116").Value = 9267.727999999999
$ws.Range("I116").Value = 9368.5
$ws.Range("K116").Value = 9368.5
$ws.Range("M116").Value = -7074.5

$ws.Range("H132").Value = 877817.5600000001
$ws.Range("I132").Value = 5523.2583
$ws.Range("J132").Value = 5384671.5
$ws.Range("K132").Value = 16569.7749
$ws.Range("L132").Value = 16154014.5
$ws.Range("M132").Value = -14039.7749
$ws.Range("N132").Value = -16159074.5

$ws.Range("H136").Value = 690346.1
$ws.Range("I136").Value = 4095.9412
$ws.Range("K136").Value = 12287.8236
$ws.Range("M136").Value = -9737.8236

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 9267.727999999999
$ws.Range("I3").Value = 9368.5
$ws.Range("K3").Value = 9368.5
$ws.Range("M3").Value = -9254.5

$ws.Range("H20").Value = 21587.934
$ws.Range("I20").Value = 7807.5
$ws.Range("K20").Value = 7807.5
$ws.Range("M20").Value = -7560.5

$ws.Range("H22").Value = 12708.889
$ws.Range("I22").Value = 12708.889
$ws.Range("K22").Value = 12708.889
$ws.Range("M22").Value = -12535.889

$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = ""
$ws.Range("N86").Value = ""

$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = ""
$ws.Range("N89").Value = ""

$ws.Range("H94").Value = 722.38464
$ws.Range("I94").Value = 707.5217
$ws.Range("K94").Value = 707.5217
$ws.Range("M94").Value = -256.5217

$ws.Range("H99").Value = 27546.666
$ws.Range("I99").Value = 53989.75
$ws.Range("J99").Value = 6392.2
$ws.Range("K99").Value = 53989.75
$ws.Range("L99").Value = 6392.2
$ws.Range("M99").Value = -52491.75
$ws.Range("N99").Value = -9388.200000000001

$ws.Range("H105").Value = 1999.3334
$ws.Range("I105").Value = 1999.3334
$ws.Range("K105").Value = 1999.3334
$ws.Range("M105").Value = -252.3334

$ws.Range("H134").Value = 21541.357
$ws.Range("I134").Value = 17173.111
$ws.Range("K134").Value = 51519.333
$ws.Range("M134").Value = -48984.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 47181668
$ws.Range("I6").Value = 47181668
$ws.Range("K6").Value = 47181668
$ws.Range("M6").Value = -47181555

$ws.Range("H16").Value = 8750.615
$ws.Range("I16").Value = 8750.615
$ws.Range("K16").Value = 8750.615
$ws.Range("M16").Value = -8463.615

$ws.Range("H19").Value = 300
$ws.Range("I19").Value = 300
$ws.Range("K19").Value = 300
$ws.Range("M19").Value = -130

$ws.Range("H24").Value = 300
$ws.Range("I24").Value = 300
$ws.Range("K24").Value = 300
$ws.Range("M24").Value = -130

$ws.Range("H31").Value = 70310.71000000001
$ws.Range("J31").Value = 29152.666
$ws.Range("L31").Value = 29152.666
$ws.Range("N31").Value = -29742.666

$ws.Range("H34").Value = 70310.71000000001
$ws.Range("J34").Value = 29152.666
$ws.Range("L34").Value = 29152.666
$ws.Range("N34").Value = -29556.666

$ws.Range("H35").Value = 1037.037
$ws.Range("I35").Value = 1037.037
$ws.Range("K35").Value = 1037.037
$ws.Range("M35").Value = -743.037

$ws.Range("H41").Value = 200024000

$ws.Range("H58").Value = 16098.074
$ws.Range("I58").Value = 5766.952
$ws.Range("K58").Value = 5766.952
$ws.Range("M58").Value = -5563.952

$ws.Range("H62").Value = 4000
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").Value = ""

$ws.Range("H65").Value = 4000
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").Value = ""

$ws.Range("H68").Value = 38495
$ws.Range("J68").Value = 38495
$ws.Range("L68").Value = 38495
$ws.Range("N68").Value = -39993

$ws.Range("H71").Value = 38495
$ws.Range("J71").Value = 38495
$ws.Range("L71").Value = 115485
$ws.Range("N71").Value = -122973

$ws.Range("H93").Value = 10000
$ws.Range("I93").Value = 10000
$ws.Range("K93").Value = 10000
$ws.Range("M93").Value = -8128

$ws.Range("H103").Value = 10500
$ws.Range("J103").Value = 10500
$ws.Range("L103").Value = 10500
$ws.Range("N103").Value = -12844

$ws.Range("H107").Value = 1439.2142
$ws.Range("I107").Value = 1695
$ws.Range("J107").Value = 799.75
$ws.Range("K107").Value = 1695
$ws.Range("L107").Value = 799.75
$ws.Range("M107").Value = 225
$ws.Range("N107").Value = -4639.75

$ws.Range("H113").Value = 8750.615
$ws.Range("I113").Value = 8750.615
$ws.Range("K113").Value = 8750.615
$ws.Range("M113").Value = -6580.615

$ws.Range("H122").Value = 2466.6667
$ws.Range("I122").Value = 2466.6667
$ws.Range("K122").Value = 7400.000100000001
$ws.Range("M122").Value = -4950.000100000001

$ws.Range("H132").Value = 109257220
$ws.Range("I132").Value = 55558804
$ws.Range("J132").Value = 205914350
$ws.Range("K132").Value = 166676412
$ws.Range("L132").Value = 617743050
$ws.Range("M132").Value = -166673882
$ws.Range("N132").Value = -617748110

$ws.Range("H134").Value = 55565084
$ws.Range("I134").Value = 4115.6924
$ws.Range("J134").Value = 200023600
$ws.Range("K134").Value = 12347.0772
$ws.Range("L134").Value = 600070800
$ws.Range("M134").Value = -9812.0772
$ws.Range("N134").Value = -600075870

$ws.Range("H136").Value = 16098.074
$ws.Range("I136").Value = 5766.952
$ws.Range("K136").Value = 17300.856
$ws.Range("M136").Value = -14750.856

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 1276.6
$ws.Range("I23").Value = 190
$ws.Range("J23").Value = 1548.25
$ws.Range("K23").Value = 570
$ws.Range("L23").Value = 4644.75
$ws.Range("M23").Value = -335
$ws.Range("N23").Value = -5114.75

$ws.Range("H46").Value = 904
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 904
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 2712
$ws.Range("M46").Value = ""
$ws.Range("N46").Value = -2894

$ws.Range("H57").Value = 2500
$ws.Range("J57").Value = 3071.4285
$ws.Range("L57").Value = 9214.2855
$ws.Range("N57").Value = -10332.2855

$ws.Range("H122").Value = 10761458
$ws.Range("J122").Value = 2833770
$ws.Range("L122").Value = 25503930
$ws.Range("N122").Value = -25508830

$ws.Range("H134").Value = 5515.089
$ws.Range("I134").Value = 1917.4375
$ws.Range("K134").Value = 5752.3125
$ws.Range("M134").Value = -682.3125

$ws.Range("H140").Value = 2189.84
$ws.Range("I140").Value = 1853.8334
$ws.Range("K140").Value = 5561.5002
$ws.Range("M140").Value = -381.5002000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 20390.25
$ws.Range("I58").Value = 20390.25
$ws.Range("K58").Value = 20390.25
$ws.Range("M58").Value = -20113.25

$ws.Range("H80").Value = 2744.5557
$ws.Range("I80").Value = 2729.1428
$ws.Range("J80").Value = 2798.5
$ws.Range("K80").Value = 2729.1428
$ws.Range("L80").Value = 2798.5
$ws.Range("M80").Value = -1731.1428
$ws.Range("N80").Value = -4794.5

$ws.Range("H83").Value = 2744.5557
$ws.Range("I83").Value = 2729.1428
$ws.Range("J83").Value = 2798.5
$ws.Range("K83").Value = 13645.714
$ws.Range("L83").Value = 13992.5
$ws.Range("M83").Value = -8653.714
$ws.Range("N83").Value = -23976.5

$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").Value = ""

$ws.Range("H122").Value = 3051.75
$ws.Range("I122").Value = 2099.5
$ws.Range("K122").Value = 6298.5
$ws.Range("M122").Value = -3848.5

$ws.Range("H132").Value = 659924.4399999999
$ws.Range("I132").Value = 4698.4707
$ws.Range("J132").Value = 1897573.5
$ws.Range("K132").Value = 14095.4121
$ws.Range("L132").Value = 5692720.5
$ws.Range("M132").Value = -11565.4121
$ws.Range("N132").Value = -5697780.5

$ws.Range("H134").Value = 156663
$ws.Range("J134").Value = 156663
$ws.Range("L134").Value = 469989
$ws.Range("N134").Value = -475059

$ws.Range("H135").Value = 245000
$ws.Range("J135").Value = 245000
$ws.Range("L135").Value = 245000
$ws.Range("N135").Value = -255140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7588.485
$ws.Range("I7").Value = 9166.056
$ws.Range("J7").Value = 5695.4
$ws.Range("K7").Value = 9166.056
$ws.Range("L7").Value = 5695.4
$ws.Range("M7").Value = -9054.056
$ws.Range("N7").Value = -5919.4

$ws.Range("H16").Value = 1422.3513
$ws.Range("I16").Value = 1327.6923
$ws.Range("J16").Value = 1646.091
$ws.Range("K16").Value = 1327.6923
$ws.Range("L16").Value = 1646.091
$ws.Range("M16").Value = -1157.6923
$ws.Range("N16").Value = -1986.091

$ws.Range("H22").Value = 1146.6666
$ws.Range("I22").Value = 750
$ws.Range("J22").Value = 1345
$ws.Range("K22").Value = 750
$ws.Range("L22").Value = 1345
$ws.Range("M22").Value = -455
$ws.Range("N22").Value = -1935

$ws.Range("H27").Value = 1146.6666
$ws.Range("I27").Value = 750
$ws.Range("J27").Value = 1345
$ws.Range("K27").Value = 750
$ws.Range("L27").Value = 1345
$ws.Range("M27").Value = -643
$ws.Range("N27").Value = -1559

$ws.Range("H38").Value = 84999.5
$ws.Range("J38").Value = 99999
$ws.Range("L38").Value = 99999
$ws.Range("N38").Value = -100819

$ws.Range("H45").Value = 21333.334
$ws.Range("I45").Value = 9000
$ws.Range("K45").Value = 9000
$ws.Range("M45").Value = -8593

$ws.Range("H46").Value = 1174.5
$ws.Range("I46").Value = 850
$ws.Range("J46").Value = 1499
$ws.Range("K46").Value = 850
$ws.Range("L46").Value = 1499
$ws.Range("M46").Value = -662
$ws.Range("N46").Value = -1875

$ws.Range("H54").Value = 78070
$ws.Range("J54").Value = 78070
$ws.Range("L54").Value = 78070
$ws.Range("N54").Value = -79358

$ws.Range("H55").Value = 595.7857
$ws.Range("J55").Value = 738.7143
$ws.Range("L55").Value = 738.7143
$ws.Range("N55").Value = -1084.7143

$ws.Range("H61").Value = 2994.6191
$ws.Range("I61").Value = 2440.9443
$ws.Range("J61").Value = 6316.6665
$ws.Range("K61").Value = 2440.9443
$ws.Range("L61").Value = 6316.6665
$ws.Range("M61").Value = -2238.9443
$ws.Range("N61").Value = -6720.6665

$ws.Range("H68").Value = 4214012.5
$ws.Range("I68").Value = 18760.666
$ws.Range("J68").Value = 8933671
$ws.Range("K68").Value = 18760.666
$ws.Range("L68").Value = 8933671
$ws.Range("M68").Value = -18011.666
$ws.Range("N68").Value = -8935169

$ws.Range("H71").Value = 4214012.5
$ws.Range("I71").Value = 18760.666
$ws.Range("J71").Value = 8933671
$ws.Range("K71").Value = 93803.33
$ws.Range("L71").Value = 44668355
$ws.Range("M71").Value = -90059.33
$ws.Range("N71").Value = -44675843

$ws.Range("H82").Value = 1388.3125
$ws.Range("I82").Value = 842.3333
$ws.Range("J82").Value = 2090.2856
$ws.Range("K82").Value = 842.3333
$ws.Range("L82").Value = 2090.2856
$ws.Range("M82").Value = -481.3333
$ws.Range("N82").Value = -2812.2856

$ws.Range("H85").Value = 1388.3125
$ws.Range("I85").Value = 842.3333
$ws.Range("J85").Value = 2090.2856
$ws.Range("K85").Value = 842.3333
$ws.Range("L85").Value = 2090.2856
$ws.Range("M85").Value = 405.6667
$ws.Range("N85").Value = -4586.2856

$ws.Range("H93").Value = 7332.65
$ws.Range("I93").Value = 7982.1113
$ws.Range("J93").Value = 1487.5
$ws.Range("K93").Value = 7982.1113
$ws.Range("L93").Value = 1487.5
$ws.Range("M93").Value = -6734.1113
$ws.Range("N93").Value = -3983.5

$ws.Range("H100").Value = 3041.3333
$ws.Range("I100").Value = 2804.625
$ws.Range("J100").Value = 3798.8
$ws.Range("K100").Value = 2804.625
$ws.Range("L100").Value = 3798.8
$ws.Range("M100").Value = -2263.625
$ws.Range("N100").Value = -4880.8

$ws.Range("H113").Value = 2994.6191
$ws.Range("I113").Value = 2440.9443
$ws.Range("J113").Value = 6316.6665
$ws.Range("K113").Value = 2440.9443
$ws.Range("L113").Value = 6316.6665
$ws.Range("M113").Value = -270.9443000000001
$ws.Range("N113").Value = -10656.6665

$ws.Range("H126").Value = 7588.485
$ws.Range("I126").Value = 9166.056
$ws.Range("J126").Value = 5695.4
$ws.Range("K126").Value = 27498.168
$ws.Range("L126").Value = 17086.2
$ws.Range("M126").Value = -25028.168
$ws.Range("N126").Value = -22026.2

$ws.Range("H132").Value = 1000875.9
$ws.Range("I132").Value = 4325.9
$ws.Range("J132").Value = 2993975.8
$ws.Range("K132").Value = 12977.7
$ws.Range("L132").Value = 8981927.399999999
$ws.Range("M132").Value = -10447.7
$ws.Range("N132").Value = -8986987.399999999

$ws.Range("H136").Value = 1283414.2
$ws.Range("I136").Value = 25146.334
$ws.Range("K136").Value = 75439.00199999999
$ws.Range("M136").Value = -72889.00199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 29999.6
$ws.Range("I51").Value = 19999.5
$ws.Range("J51").Value = 70000
$ws.Range("K51").Value = 19999.5
$ws.Range("L51").Value = 70000
$ws.Range("M51").Value = -19489.5
$ws.Range("N51").Value = -71020

$ws.Range("H92").Value = 2366.6667
$ws.Range("I92").Value = 2000
$ws.Range("J92").Value = 2550
$ws.Range("K92").Value = 2000
$ws.Range("L92").Value = 2550
$ws.Range("M92").Value = 496
$ws.Range("N92").Value = -7542

$ws.Range("H96").Value = 1931
$ws.Range("I96").Value = 1193
$ws.Range("J96").Value = 2300
$ws.Range("K96").Value = 1193
$ws.Range("L96").Value = 2300
$ws.Range("M96").Value = 180
$ws.Range("N96").Value = -5046

$ws.Range("H122").Value = 2973.75
$ws.Range("J122").Value = 3750
$ws.Range("L122").Value = 11250
$ws.Range("N122").Value = -16150

$ws.Range("H132").Value = 2935135.8
$ws.Range("I132").Value = 7132.6665
$ws.Range("J132").Value = 7327140
$ws.Range("K132").Value = 21397.9995
$ws.Range("L132").Value = 21981420
$ws.Range("M132").Value = -18867.9995
$ws.Range("N132").Value = -21986480

$ws.Range("H136").Value = 345092.25
$ws.Range("I136").Value = 3142.5293
$ws.Range("J136").Value = 760316.9399999999
$ws.Range("K136").Value = 9427.5879
$ws.Range("L136").Value = 2280950.82
$ws.Range("M136").Value = -6877.5879
$ws.Range("N136").Value = -2286050.82
